$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.95
$ws.Range("I3").Value = 4.5
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.5
$ws.Range("U3").Value = 4
$ws.Range("V3").Value = 1.23
$ws.Range("AG3").Value = 19
$ws.Range("AN3").Value = 9.5
$ws.Range("AO3").Value = 21
$ws.Range("AP3").Value = 17
$ws.Range("AQ3").Value = 51

# Row 4 updates
$ws.Range("S4").Value = 2.2
$ws.Range("T4").Value = 1.65
